$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 99, pushing the existing
# rows 99-108 down to 100-109 (keeping all their data/formatting intact).
# Restrict the insert to the used columns (A:R) so we don't stamp
# formatting/styles across the whole 16384-column row.
$ws.Range("A99:R99").Insert(-4121)  # xlShiftDown

# New row 99 gets the same layout/style as the row below it (old row 99,
# now row 100) so the date column keeps its date/time number format.
$ws.Range("A100:R100").Copy()
$ws.Range("A99:R99").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new weekly price observation in row 99.
$ws.Cells.Item(99, 1).Value = 7
$ws.Cells.Item(99, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(99, 3).Value = "Ñuble"
$ws.Cells.Item(99, 4).Value = 45131
$ws.Cells.Item(99, 5).Value = 16
$ws.Cells.Item(99, 6).Value = 100112013
$ws.Cells.Item(99, 7).Value = "Alcachofa"
$ws.Cells.Item(99, 8).Value = "Madrigal"
$ws.Cells.Item(99, 9).Value = "Primera"
$ws.Cells.Item(99, 10).Value = 50
$ws.Cells.Item(99, 11).Value = 15000
$ws.Cells.Item(99, 12).Value = 15000
$ws.Cells.Item(99, 13).Value = 15000
$ws.Cells.Item(99, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(99, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(99, 16).Value = 375
$ws.Cells.Item(99, 17).Value = 40
$ws.Cells.Item(99, 18).Value = "Hortaliza"
